$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- Row 7 (Step "1.") -------------------------------------------------
# Opis kroku testowego (C7)
$ws.Range("C7").Value = "Wejdź na główną stronę sklepu, usuń pliki cookies i zweryfikuj tytuł strony. Dodatkowo wykonaj zrzut ekranu."

# Dane testowe (D7) - text unchanged, just re-assert it
$ws.Range("D7").Value = "Adres głównej strony sklepu: http://koszulkifootball.sellingo.pl/" + $nl + "Poprawny tytuł strony to: 'Koszulkifootball.sellingo.pl'"

# Oczekiwany / Aktualny rezultat (E7 / F7)
$row7Result = "Strona główna sklepu wyświetlona" + $nl + "Tytuł strony: ''Koszulkifootball.sellingo.pl''" + $nl + "Pliki Cookies usunięte" + $nl + "Zrzut ekranu zapisany do pliku .jpg"
$ws.Range("E7").Value = $row7Result
$ws.Range("F7").Value = $row7Result

# --- Row 9 (Step "3.") -------------------------------------------------
# Opis kroku testowego (C9)
$ws.Range("C9").Value = "Kliknij w zakładkę `"Kontakt' i zweryfikuj tytuł  strony."

# Dane testowe (D9) - text unchanged, just re-assert it
$ws.Range("D9").Value = "Poprawny tytuł strony to: 'Kontakt'"

# Oczekiwany / Aktualny rezultat (E9 / F9) - same text plus trailing newline
$row9Result = "Użytkownik zostaje przeniesiony na podstronę 'Kontakt'. Tytuł strony: 'Kontakt'" + $nl
$ws.Range("E9").Value = $row9Result
$ws.Range("F9").Value = $row9Result

# --- Restore the last active selection ---------------------------------
$ws.Activate()
$ws.Range("F10").Select()
